# Apply updated "K" column (column G) values to the save_data sheet.
# The commit regenerates the K (strike count) column from newly recomputed
# s_vals, replacing the previous (stale "Strike#"-based) numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 4
    3  = 6
    4  = 4
    5  = 1
    6  = 5
    7  = 4
    8  = 3
    9  = 8
    10 = 4
    11 = 5
    12 = 9
    13 = 3
    14 = 4
    15 = 8
    16 = 6
    17 = 6
    18 = 1
    19 = 1
    20 = 1
    21 = 8
    22 = 5
    23 = 1
    24 = 4
    25 = 2
    26 = 3
    27 = 6
    28 = 8
    29 = 7
    30 = 2
    31 = 5
    32 = 1
    33 = 1
    34 = 3
    35 = 3
    36 = 3
    37 = 7
    38 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
